$d = $word.ActiveDocument
$sec = $d.Sections.First

# Header 1 (the BTEC logo, Y-less local drawing) currently has its picture
# labelled "image1.jpg" on both the wp:docPr and pic:cNvPr nodes; the
# commit renames that picture to "image2.jpg".
$hdr = $sec.Headers.Item(2)
if ($hdr.Exists) {
    $shapes = $hdr.Range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
            $shape.Name = "image2.jpg"
        }
    }
}

# Footer 1 and Footer 2 both carry the Pearson Edexcel logo, currently
# named "image2.png"; the commit renames both occurrences to "image1.png".
for ($f = 1; $f -le 2; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shape = $shapes.Item($i)
            if ($shape.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shape.Name = "image1.png"
            }
        }
    }
}
